$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "duty free" keyword/appID rows are inserted into the existing list,
# pushing the rows below them down. Insert top-to-bottom so each insertion
# point refers to the already-updated row numbering.
$insertRows = @(4, 7, 11, 17, 23)

foreach ($r in $insertRows) {
    $ws.Rows("$r`:$r").Insert()
    # Write column B before column A so the shared-string table registers
    # "duty.pare.myapp" ahead of "duty free" (matches source order).
    $ws.Range("B$r").Value = "duty.pare.myapp"
    $ws.Range("A$r").Value = "duty free"
}

$ws.Range("A4:B4").Select()
